# Add five tests for login with different valid usernames.
# This re-orders the "locked_out_user" row to come right after
# "standard_user" (row 3), shifts the remaining existing rows down
# by one, and appends a brand-new row for "visual_user".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current rows 3-6 (A/C/D columns; B is unused below row 2)
# before we start overwriting them.
$row3 = @($ws.Range("A3").Value2, $ws.Range("C3").Value2, $ws.Range("D3").Value2)
$row4 = @($ws.Range("A4").Value2, $ws.Range("C4").Value2, $ws.Range("D4").Value2)
$row5 = @($ws.Range("A5").Value2, $ws.Range("C5").Value2, $ws.Range("D5").Value2)
$row6 = @($ws.Range("A6").Value2, $ws.Range("C6").Value2, $ws.Range("D6").Value2)

# New row order: locked_out_user moves up to row 3, the previous
# row3/4/5 (problem_user, performance_glitch_user, error_user) shift
# down to rows 4/5/6.
$ws.Range("A3").Value = $row6[0]
$ws.Range("C3").Value = $row6[1]
$ws.Range("D3").Value = $row6[2]

$ws.Range("A4").Value = $row3[0]
$ws.Range("C4").Value = $row3[1]
$ws.Range("D4").Value = $row3[2]

$ws.Range("A5").Value = $row4[0]
$ws.Range("C5").Value = $row4[1]
$ws.Range("D5").Value = $row4[2]

$ws.Range("A6").Value = $row5[0]
$ws.Range("C6").Value = $row5[1]
$ws.Range("D6").Value = $row5[2]

# New row 7 for the added test case.
$ws.Range("A7").Value = "visual_user"
$ws.Range("C7").Value = "visualuser"
$ws.Range("D7").Value = "saucelabs"

# Update the active selection to match the saved view state.
$ws.Range("C12").Select()
